$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new 2023 column (T), matching the formatting of the preceding column (S)
$ws.Range("S4:S5").Copy()
$ws.Range("T4:T5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("T4").Value = 2023
$ws.Range("T5").Value = 40

# Narrow columns A:C slightly (39 chars -> ~36.57 chars)
$ws.Range("A1:C1").ColumnWidth = 35.65

# Reset the view back to A1 (clears the stored topLeftCell="C1" / U4 selection)
$ws.Range("A1").Select() | Out-Null
